# Weekly CompStat report refresh: new crime data collected.
# Updates header volume/number + reporting week dates, and refreshes the
# crime-complaint statistics table (rows 15-30) with new figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header text: "Volume 31   Number  7" -> "...  8" and the reporting
# week dates "2/12/2024 .. 2/18/2024" -> "2/19/2024 .. 2/25/2024"
# ---------------------------------------------------------------------
$ws.Cells.Item(8, 1).Value = "Volume 31   Number  8"
$ws.Cells.Item(9, 3).Value = "Report Covering the Week  2/19/2024  Through  2/25/2024"

# ---------------------------------------------------------------------
# Helper: write a plain numeric value into a cell.
# ---------------------------------------------------------------------
function Set-Num($addr, $value) {
    $ws.Range($addr).Value = $value
}

# Helper: force a cell to hold the text "0" (shared text used for an
# empty/zero count column) using the same visual style as other
# existing text cells in the table (style carried by C15).
function Set-TextZero($addr) {
    $ws.Range($addr).Value = "'0"
    $ws.Range("C15").Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
}

# Helper: force a cell to hold the text "***.*" (used when a percent
# change is undefined), matching the style already used by N22.
function Set-TextStar($addr) {
    $ws.Range($addr).Value = "'***.*"
    $ws.Range("N22").Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
}

# Helper: turn a text "0" cell back into a plain number cell, matching
# the numeric style already used by D16.
function Set-NumFromText($addr, $value) {
    $ws.Range("D16").Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
    $ws.Range($addr).Value = $value
}

$ws.Application.CutCopyMode = 0

# ---------------------------------------------------------------------
# Row 15 - Rape
# ---------------------------------------------------------------------
Set-TextZero "D15"
Set-TextStar "E15"

# ---------------------------------------------------------------------
# Row 16 - Robbery
# ---------------------------------------------------------------------
Set-Num "C16" 2
Set-Num "D16" 2
Set-Num "E16" 0
Set-Num "G16" 13
Set-Num "H16" -7.692307692307
Set-Num "I16" 24
Set-Num "J16" 28
Set-Num "K16" -14.285714285714
Set-Num "L16" -7.692307692307
Set-Num "M16" -14.285714285714
Set-Num "N16" -82.481751824817

# ---------------------------------------------------------------------
# Row 17 - Fel. Assault
# ---------------------------------------------------------------------
Set-TextZero "C17"
Set-Num "D17" 8
Set-Num "E17" -100
Set-Num "F17" 13
Set-Num "H17" -7.142857142857
Set-Num "J17" 37
Set-Num "K17" -37.837837837837
Set-Num "L17" -34.285714285714
Set-Num "M17" 53.333333333333
Set-Num "N17" -32.352941176470

# ---------------------------------------------------------------------
# Row 18 - Burglary
# ---------------------------------------------------------------------
Set-Num "C18" 2
Set-Num "D18" 5
Set-Num "E18" -60
Set-Num "F18" 11
Set-Num "H18" -42.105263157894
Set-Num "I18" 35
Set-Num "J18" 38
Set-Num "K18" -7.894736842105
Set-Num "L18" -5.405405405405
Set-Num "M18" -39.655172413793
Set-Num "N18" -88.135593220339

# ---------------------------------------------------------------------
# Row 19 - Gr. Larceny
# ---------------------------------------------------------------------
Set-Num "C19" 19
Set-Num "D19" 21
Set-Num "E19" -9.523809523809
Set-Num "F19" 62
Set-Num "G19" 78
Set-Num "H19" -20.512820512820
Set-Num "I19" 147
Set-Num "J19" 141
Set-Num "K19" 4.255319148936
Set-Num "L19" 6.521739130434
Set-Num "M19" -30.985915492957
Set-Num "N19" -63.157894736842

# ---------------------------------------------------------------------
# Row 20 - G.L.A.
# ---------------------------------------------------------------------
Set-Num "F20" 2
Set-Num "G20" 7
Set-Num "H20" -71.428571428571
Set-Num "J20" 8
Set-Num "K20" -50
Set-Num "L20" -55.555555555555
Set-Num "N20" -98.076923076923

# ---------------------------------------------------------------------
# Row 21 - TOTAL
# ---------------------------------------------------------------------
Set-Num "D21" 39
Set-Num "E21" -41.025641025641
Set-Num "F21" 100
Set-Num "G21" 132
Set-Num "H21" -24.242424242424
Set-Num "I21" 233
Set-Num "J21" 253
Set-Num "K21" -7.905138339920
Set-Num "L21" -6.425702811244
Set-Num "M21" -27.1875
Set-Num "N21" -78.365831012070

# ---------------------------------------------------------------------
# Row 22 - Transit
# ---------------------------------------------------------------------
Set-NumFromText "C22" 2
Set-Num "E22" -50
Set-Num "G22" 10
Set-Num "H22" -40
Set-Num "I22" 16
Set-Num "J22" 16
Set-Num "K22" 0
Set-Num "L22" -15.789473684210
Set-Num "M22" 45.454545454545

# ---------------------------------------------------------------------
# Row 23 - Housing
# ---------------------------------------------------------------------
Set-Num "G23" 2
Set-Num "J23" 3
Set-Num "K23" -66.666666666666

# ---------------------------------------------------------------------
# Row 24 - Petit Larceny
# ---------------------------------------------------------------------
Set-Num "C24" 69
Set-Num "D24" 32
Set-Num "E24" 115.625
Set-Num "F24" 223
Set-Num "G24" 153
Set-Num "H24" 45.751633986928
Set-Num "I24" 461
Set-Num "J24" 304
Set-Num "K24" 51.644736842105
Set-Num "L24" 31.339031339031
Set-Num "M24" 68.248175182481

# ---------------------------------------------------------------------
# Row 25 - Misd. Assault
# ---------------------------------------------------------------------
Set-Num "C25" 12
Set-Num "D25" 7
Set-Num "E25" 71.428571428571
Set-Num "F25" 50
Set-Num "G25" 31
Set-Num "H25" 61.290322580645
Set-Num "I25" 80
Set-Num "J25" 63
Set-Num "K25" 26.984126984127
Set-Num "L25" 35.593220338983
Set-Num "M25" 56.862745098039

# ---------------------------------------------------------------------
# Row 26 - UCR Rape*
# ---------------------------------------------------------------------
Set-TextZero "D26"
Set-TextStar "E26"
Set-Num "G26" 3

# ---------------------------------------------------------------------
# Row 27 - Other Sex Crimes
# ---------------------------------------------------------------------
Set-Num "C27" 1
Set-TextZero "D27"
Set-TextStar "E27"
Set-Num "G27" 4
Set-Num "H27" 25
Set-Num "I27" 14
Set-Num "K27" 40
Set-Num "L27" 0

# ---------------------------------------------------------------------
# Row 30 - Hate Crimes
# ---------------------------------------------------------------------
Set-Num "F30" 3
Set-TextZero "G30"
Set-TextStar "H30"
Set-Num "I30" 4
Set-Num "K30" 300
Set-Num "L30" 100

$ws.Application.CutCopyMode = 0
